$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 24305.172
$ws.Range("J51").Value = 25411.969
$ws.Range("L51").Value = 25411.969
$ws.Range("N51").Value = -26379.969
$ws.Range("H63").Value = 113999
$ws.Range("J63").Value = 113999
$ws.Range("L63").Value = 113999
$ws.Range("N63").Value = -115247
$ws.Range("H66").Value = 113999
$ws.Range("J66").Value = 113999
$ws.Range("L66").Value = 341997
$ws.Range("N66").Value = -348237
$ws.Range("H100").Value = 6583.3335
$ws.Range("I100").Value = 7375
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 7375
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -6834
$ws.Range("N100").Value = -6082
$ws.Range("H113").Value = 71435910
$ws.Range("I113").Value = 25007850
$ws.Range("J113").Value = 133340000
$ws.Range("K113").Value = 25007850
$ws.Range("L113").Value = 133340000
$ws.Range("M113").Value = -25004596
$ws.Range("N113").Value = -133346508
$ws.Range("H138").Value = 3514.4722
$ws.Range("J138").Value = 3333.3823
$ws.Range("L138").Value = 10000.1469
$ws.Range("N138").Value = -20280.1469

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 333.33334
$ws.Range("J6").Value = 200
$ws.Range("L6").Value = 200
$ws.Range("N6").Value = -546
$ws.Range("H32").Value = 6677947.5
$ws.Range("I32").Value = 7048113.5
$ws.Range("K32").Value = 7048113.5
$ws.Range("M32").Value = -7047826.5
$ws.Range("H61").Value = 17284718
$ws.Range("I61").Value = 35719588
$ws.Range("J61").Value = 78837.87
$ws.Range("K61").Value = 35719588
$ws.Range("L61").Value = 78837.87
$ws.Range("M61").Value = -35719376
$ws.Range("N61").Value = -79261.87
$ws.Range("H136").Value = 17284718
$ws.Range("I136").Value = 35719588
$ws.Range("J136").Value = 78837.87
$ws.Range("K136").Value = 107158764
$ws.Range("L136").Value = 236513.61
$ws.Range("M136").Value = -107156214
$ws.Range("N136").Value = -241613.61

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 8000
$ws.Range("J63").Value = 8000
$ws.Range("L63").Value = 8000
$ws.Range("N63").Value = -9372
$ws.Range("H66").Value = 8000
$ws.Range("J66").Value = 8000
$ws.Range("L66").Value = 24000
$ws.Range("N66").Value = -30864
$ws.Range("H94").Value = 556.8929000000001
$ws.Range("I94").Value = 523.28
$ws.Range("K94").Value = 523.28
$ws.Range("M94").Value = -72.27999999999997
$ws.Range("H140").Value = 249999.33
$ws.Range("J140").Value = 249999.33
$ws.Range("L140").Value = 249999.33
$ws.Range("N140").Value = -260359.33

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H64").Value = 29663.334
$ws.Range("J64").Value = 29663.334
$ws.Range("L64").Value = 29663.334
$ws.Range("N64").Value = -30159.334
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H67").Value = 29663.334
$ws.Range("J67").Value = 29663.334
$ws.Range("L67").Value = 29663.334
$ws.Range("N67").Value = -31379.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 30857.77
$ws.Range("J2").Value = 66773.836
$ws.Range("L2").Value = 400643.0159999999
$ws.Range("N2").Value = -400869.0159999999
$ws.Range("H21").Value = 1966.6666
$ws.Range("J21").Value = 900
$ws.Range("L21").Value = 2700
$ws.Range("N21").Value = -3046
$ws.Range("H23").Value = 330.55554
$ws.Range("I23").Value = 306.25
$ws.Range("J23").Value = 350
$ws.Range("K23").Value = 918.75
$ws.Range("L23").Value = 1050
$ws.Range("M23").Value = -683.75
$ws.Range("N23").Value = -1520
$ws.Range("H68").Value = 1675.75
$ws.Range("I68").Value = 1723.8334
$ws.Range("J68").Value = 1664.6538
$ws.Range("K68").Value = 5171.5002
$ws.Range("L68").Value = 4993.9614
$ws.Range("M68").Value = -4360.5002
$ws.Range("N68").Value = -6615.9614
$ws.Range("H71").Value = 1675.75
$ws.Range("I71").Value = 1723.8334
$ws.Range("J71").Value = 1664.6538
$ws.Range("K71").Value = 15514.5006
$ws.Range("L71").Value = 14981.8842
$ws.Range("M71").Value = -11458.5006
$ws.Range("N71").Value = -23093.8842
$ws.Range("H86").Value = 349.25
$ws.Range("J86").Value = 299
$ws.Range("L86").Value = 897
$ws.Range("N86").Value = -3269
$ws.Range("H89").Value = 349.25
$ws.Range("J89").Value = 299
$ws.Range("L89").Value = 2691
$ws.Range("N89").Value = -14547
$ws.Range("H113").Value = 1826.8235
$ws.Range("I113").Value = 1428.3334
$ws.Range("J113").Value = 1912.2142
$ws.Range("K113").Value = 4285.0002
$ws.Range("L113").Value = 5736.642599999999
$ws.Range("M113").Value = -2115.0002
$ws.Range("N113").Value = -10076.6426
$ws.Range("H129").Value = 65980
$ws.Range("I129").Value = 997.1429000000001
$ws.Range("J129").Value = 122840
$ws.Range("K129").Value = 2991.4287
$ws.Range("L129").Value = 368520
$ws.Range("M129").Value = 2008.5713
$ws.Range("N129").Value = -378520
$ws.Range("H137").Value = 5872.8887
$ws.Range("I137").Value = 4002.4666
$ws.Range("J137").Value = 15225
$ws.Range("K137").Value = 12007.3998
$ws.Range("L137").Value = 45675
$ws.Range("M137").Value = -6907.399800000001
$ws.Range("N137").Value = -55875
$ws.Range("H141").Value = 95955.06
$ws.Range("I141").Value = 106821.27
$ws.Range("K141").Value = 320463.81
$ws.Range("M141").Value = -315283.81

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 100
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H80").Value = 3652.3076
$ws.Range("I80").Value = 3332.6667
$ws.Range("J80").Value = 4371.5
$ws.Range("K80").Value = 3332.6667
$ws.Range("L80").Value = 4371.5
$ws.Range("M80").Value = -2334.6667
$ws.Range("N80").Value = -6367.5
$ws.Range("H83").Value = 3652.3076
$ws.Range("I83").Value = 3332.6667
$ws.Range("J83").Value = 4371.5
$ws.Range("K83").Value = 16663.3335
$ws.Range("L83").Value = 21857.5
$ws.Range("M83").Value = -11671.3335
$ws.Range("N83").Value = -31841.5
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H122").Value = 2219.7778
$ws.Range("I122").Value = 2244.75
$ws.Range("J122").Value = 2199.8
$ws.Range("K122").Value = 6734.25
$ws.Range("L122").Value = 6599.400000000001
$ws.Range("M122").Value = -4284.25
$ws.Range("N122").Value = -11499.4
$ws.Range("H132").Value = 57145070
$ws.Range("I132").Value = 60608296
$ws.Range("K132").Value = 181824888
$ws.Range("M132").Value = -181822358

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2249.158
$ws.Range("J68").Value = 3283
$ws.Range("L68").Value = 3283
$ws.Range("N68").Value = -4781
$ws.Range("H71").Value = 2249.158
$ws.Range("J71").Value = 3283
$ws.Range("L71").Value = 16415
$ws.Range("N71").Value = -23903
$ws.Range("H119").Value = 49990
$ws.Range("J119").Value = 49990
$ws.Range("L119").Value = 49990
$ws.Range("N119").Value = -59666
$ws.Range("H132").Value = 4835441
$ws.Range("I132").Value = 693705
$ws.Range("J132").Value = 14302266
$ws.Range("K132").Value = 2081115
$ws.Range("L132").Value = 42906798
$ws.Range("M132").Value = -2078585
$ws.Range("N132").Value = -42911858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 390214.62
$ws.Range("I132").Value = 4643.1113
$ws.Range("J132").Value = 1257750.5
$ws.Range("K132").Value = 13929.3339
$ws.Range("L132").Value = 3773251.5
$ws.Range("M132").Value = -11399.3339
$ws.Range("N132").Value = -3778311.5
